{"js": "// The document has four paragraphs that each spell out, run by run, a\n// sentence about the Perseids (\"Dates de la campanya 2018 en qu\u00e8 usem\n// la constel\u00b7laci\u00f3 Perseus 30 d'octubre al novembre 8 i 29 de novembre\n// de desembre 8\"). They all get collapsed to a single plain run with\n// the new, translated Orion campaign dates.\nconst newText =\n  \"Dates de la campanya Orion: 16-25 de gener, 14-23 de febrer, del 14 al 24 de mar\u00e7\";\n\nconst body = context.document.body;\n\n// Case-sensitive search so the unrelated lower-case mention (\"...sempre\n// dins de les dates de la campanya.\") elsewhere in the document is left\n// untouched.\nconst results = body.search(\"Dates de la campanya\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\n// Resolve the owning paragraph for every hit up front (paragraphs shift\n// as earlier ones are rewritten, so grab them all before mutating).\nconst paragraphs = [];\nfor (let i = 0; i < results.items.length; i++) {\n  paragraphs.push(results.items[i].paragraphs.getFirst());\n}\nawait context.sync();\n\n// Clear each paragraph's whole range first: deleting its (many,\n// variously-formatted) runs so the text inserted afterwards starts a\n// brand-new run with no inherited character formatting, matching the\n// target OOXML (a bare <w:r><w:t>\u2026</w:t></w:r>).\nconst ranges = paragraphs.map((p) => p.getRange());\nfor (const range of ranges) {\n  range.clear();\n}\nawait context.sync();\n\nfor (const range of ranges) {\n  range.insertText(newText, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# The document has four paragraphs that each spell out, run by run, a\n# sentence about the Perseids (\"Dates de la campanya 2018 en qu\u00e8 usem\n# la constel\u00b7laci\u00f3 Perseus 30 d'octubre al novembre 8 i 29 de novembre\n# de desembre 8\"). They all get collapsed to a single plain run with\n# the new, translated Orion campaign dates.\n$newText = \"Dates de la campanya Orion: 16-25 de gener, 14-23 de febrer, del 14 al 24 de mar\u00e7\"\n\n$d = $word.ActiveDocument\n\n# Locate every target paragraph up front (by a fragment - \"30 d'octubre\"\n# - unique to these four paragraphs) before mutating anything, so a live\n# $d.Paragraphs enumerator isn't disturbed mid-loop. This also avoids the\n# one unrelated, lower-case mention (\"...sempre dins de les dates de la\n# campanya.\") elsewhere in the document.\n$targets = @()\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains(\"30 d'octubre\")) {\n        $targets += $p\n    }\n}\n\nforeach ($p in $targets) {\n    $r = $p.Range\n    # Exclude the trailing paragraph mark, then delete the whole run\n    # sequence and insert fresh text - this starts a brand-new run with\n    # no inherited character formatting, matching the target OOXML (a\n    # bare <w:r><w:t>...</w:t></w:r>) instead of carrying over whichever\n    # run's rPr a plain Range.Text assignment would keep.\n    $r.MoveEnd(1, -1)\n    $r.Delete()\n    $r.InsertAfter($newText)\n}\n"}
